$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Range("C16").Value = "1007208457"
$ws.Range("D16").Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Range("E16").Value = "2007"
$ws.Range("F16").Value = 24292
$ws.Range("C17").Value = "1007208457"
$ws.Range("D17").Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Range("E17").Value = "2006"
$ws.Range("F17").Value = 33125
$ws.Range("C18").Value = "1007208457"
$ws.Range("D18").Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Range("E18").Value = "2005"
$ws.Range("F18").Value = 33125
$ws.Range("C19").Value = "1007208457"
$ws.Range("D19").Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Range("E19").Value = "2004"
$ws.Range("F19").Value = 33125
$ws.Range("C20").Value = "1007208457"
$ws.Range("D20").Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Range("E20").Value = "2003"
$ws.Range("F20").Value = 33125
$ws.Range("C21").Value = "1007208457"
$ws.Range("D21").Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Range("E21").Value = "2002"
$ws.Range("F21").Value = 33125
$ws.Range("C22").Value = "1007208457"
$ws.Range("D22").Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Range("E22").Value = "2001"
$ws.Range("F22").Value = 33125
$ws.Range("C23").Value = "1007208457"
$ws.Range("D23").Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Range("E23").Value = "1912"
$ws.Range("F23").Value = 33125
$ws.Range("C24").Value = "1007208457"
$ws.Range("D24").Value = "MIGUEL ANTONIO ELLES BANGUERA"
$ws.Range("E24").Value = "1911"
$ws.Range("F24").Value = 33125
$ws.Range("C25").Value = "1049927922"
$ws.Range("D25").Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Range("E25").Value = "2007"
$ws.Range("F25").Value = 24292
$ws.Range("C26").Value = "1049927922"
$ws.Range("D26").Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Range("E26").Value = "2006"
$ws.Range("F26").Value = 33125
$ws.Range("C27").Value = "1049927922"
$ws.Range("D27").Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Range("E27").Value = "2005"
$ws.Range("F27").Value = 33125
$ws.Range("C28").Value = "1049927922"
$ws.Range("D28").Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Range("E28").Value = "2004"
$ws.Range("F28").Value = 33125
$ws.Range("C29").Value = "1049927922"
$ws.Range("D29").Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Range("E29").Value = "2003"
$ws.Range("F29").Value = 33125
$ws.Range("C30").Value = "1049927922"
$ws.Range("D30").Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Range("E30").Value = "2002"
$ws.Range("F30").Value = 33125
$ws.Range("C31").Value = "1049927922"
$ws.Range("D31").Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Range("E31").Value = "2001"
$ws.Range("F31").Value = 33125
$ws.Range("C32").Value = "1049927922"
$ws.Range("D32").Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Range("E32").Value = "1912"
$ws.Range("F32").Value = 33125
$ws.Range("C33").Value = "1049927922"
$ws.Range("D33").Value = "WILMER DE JESUS MARTINEZ DE LOS REYES"
$ws.Range("E33").Value = "1911"
$ws.Range("F33").Value = 33125
